# Update the ByName Command Sequence diagram.
#
# Re-aligns several shapes on slide 1 of the diagram and drops the stale
# start-connection on "Straight Connector 90" (id 91), which used to snap
# to "Rectangle 62" (id 90).
#
# NOTE on the literal Top values below: PowerPoint COM measures Shape.Top
# in points, but the file stores EMU (1 pt = 12700 EMU). The point values
# here were chosen (offline) so that converting back to EMU lands exactly
# on the target offsets used by the authored diagram.

function Get-ShapeById {
    param($Slide, [int]$Id)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $candidate = $Slide.Shapes.Item($i)
        if ($candidate.Id -eq $Id) {
            return $candidate
        }
    }
    throw "Shape with id $Id not found"
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Rectangle 65 (id 81) - outer "Logic" frame moves up slightly. (y: 115757 -> 76200)
$sh81 = Get-ShapeById $s 81
$sh81.Top = 6.00003937007874

# Rectangle 62 (id 90) - ":ByNameCommand" box moves down. (y: 4141907 -> 4253826)
$sh90 = Get-ShapeById $s 90
$sh90.Top = 334.94696042999504

# Straight Connector 90 (id 91) - drop its stale start-connection to shape
# 90 and update its new start position. (y: 4512091 -> 4532443)
$sh91 = Get-ShapeById $s 91
$sh91.ConnectorFormat.BeginDisconnect()
$sh91.Top = 356.88531496062996

# Straight Arrow Connector 91 (id 92) moves down. (y: 4333380 -> 4445299)
$sh92 = Get-ShapeById $s 92
$sh92.Top = 350.0235901359498

# Straight Arrow Connector 93 (id 94) moves down. (y: 5069681 -> 5181600)
$sh94 = Get-ShapeById $s 94
$sh94.Top = 408.00003052718995

# Rectangle 101 (id 102) moves down. (y: 4519825 -> 4631744)
$sh102 = Get-ShapeById $s 102
$sh102.Top = 364.70428464874504
